$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish row 10: add PriceChange (X10) and UpDown (Y10) ---
$ws.Cells.Item(10, 24).Value = -0.16000300000000323
$ws.Cells.Item(10, 25).Value = "Down"

# --- New row 11 (trade result from the repeater) ---
$ws.Cells.Item(11, 1).Value = 42654.894490740742
$ws.Cells.Item(11, 2).Value = 16
$ws.Cells.Item(11, 3).Value = "Strong Buy"
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = "Random"
$ws.Cells.Item(11, 17).Value = 28.689659976213832
$ws.Cells.Item(11, 18).Value = 0.84
$ws.Cells.Item(11, 19).Value = -0.0125
$ws.Cells.Item(11, 19).NumberFormat = "0.00%"
$ws.Cells.Item(11, 20).Value = -0.0261
$ws.Cells.Item(11, 20).NumberFormat = "0.00%"
$ws.Cells.Item(11, 21).Value = 14.56
$ws.Cells.Item(11, 22).Value = "N/A"
$ws.Cells.Item(11, 23).Value = 0
